$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New CRM-accuracy rows for the batch opened on 2021-03-14 (dmb 0314)
$newRows = @(
    @{ Row = 20; A = 20210314; B = 2205.7312106719601; C = 2224.4699999999998; E = 180 },
    @{ Row = 21; A = 20210314; B = 2194.00552102248;   C = 2224.4699999999998; E = 180 },
    @{ Row = 22; A = 20210314; B = 2192.0447529285698; C = 2224.4699999999998; E = 180 },
    @{ Row = 23; A = 20210314; B = 2196.3789375564002; C = 2224.4699999999998; E = 180 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Range("D$row").Formula = "=100*(B$row-C$row)/C$row"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = "CRM opened 20210314"
}

# Match the selection left behind after entering the new batch's data
$ws.Range("E21:F23").Select() | Out-Null
